# Update the "card" column (A) on sheet Card11 for the tone-range rows
# that belong to card group 11 — rows that were mistakenly left as "2"
# become "11", matching rows 2 and 9 which already read "11".
#
# The source data is stored as text (pandas export artifact), so a plain
# Value assignment of a numeric-looking string like "11" would silently
# be re-typed as a Number by Excel's smart entry. Prefixing with an
# apostrophe forces a literal text entry (quote-prefixed), and resetting
# the cell style back to "Normal" afterwards clears the quote-prefix
# formatting bit so the cell ends up as plain text with the default
# style, exactly like the surrounding text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card11")

$rows = 3,4,5,6,7,8,10,11,12,13
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "'11"
    $cell.Style = "Normal"
}
